$d = $word.ActiveDocument

# Mapping of old text -> new text (each old value is unique in the document)
$replacements = @(
    @("2024-12-09 Monday", "2024-12-10 Tuesday"),
    @("21×26=", "54×48="),
    @("65×99=", "54×45="),
    @("30×26=", "36×96="),
    @("21×93=", "22×22="),
    @("73×80=", "90×19="),
    @("34×22=", "31×34="),
    @("95×34=", "48×82="),
    @("48×59=", "97×11="),
    @("66×42=", "75×59="),
    @("22×43=", "50×53="),
    @("39×99=", "41×98="),
    @("92×33=", "13×80="),
    @("95×47=", "40×44="),
    @("28×63=", "31×86="),
    @("64×45=", "90×39="),
    @("65×39=", "61×37="),
    @("78×39=", "33×85="),
    @("90×81=", "20×18="),
    @("87×40=", "17×67="),
    @("84×13=", "59×82="),
    @("13×42=", "82×39="),
    @("93×52=", "20×96="),
    @("20×69=", "30×43="),
    @("90×55=", "52×62="),
    @("54×87=", "92×13=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
